$d = $word.ActiveDocument

# Locate the "Programa" body paragraph holding the run-on syllabus text
# (rather than a hard-coded paragraph index, so the script is resilient
# to unrelated structural shifts elsewhere in the document).
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("1. Introdução a Pesquisa Operacional1.1.")) {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not locate the 'Programa' syllabus paragraph"
}

# [char]11 is a vertical-tab, which Word's object model treats as a
# manual line break (<w:br/>) inside a run when written via Range.Text -
# splitting the single run-on sentence into one <w:t> per topic line,
# each separated by a <w:br/>, without starting new paragraphs.
$lb = [char]11
$lines = @(
    "1. Introdução a Pesquisa Operacional",
    "1.1. Conceitos de Pesquisa Operacional;",
    "1.2. Modelagem;",
    "1.3. Estrutura dos Modelos Matemáticos;",
    "1.4. Técnicas matemáticas em Pesquisa Operacional;",
    "1.2. Fases de Um Estudo em Pesquisa Operacional",
    "2. Programação Linear",
    "2.1. Definição",
    "2.2. Formulação de Modelos",
    "2.3. Resolução Gráfica;",
    "3. Método Simplex",
    "3.1. Desenvolvimento do Método Simplex;",
    "3.2. Procedimento do Método Simplex;",
    "4. Introdução aos Grafos e à Otimização em Rede",
    "4.1. Conceitos Básicos em Teoria dos Grafos",
    "4.2. Problemas de Fluxo Máximo;",
    "4.3. Problemas de Caminho Mínimo",
    "5. Estudo de Casos em Programação Linear",
    "5.1. Modelo de Transporte Simples",
    "5.2. Modelo da Designação.",
    "6. Introdução a Teoria das Filas",
    "6.1. Conceitos da Teoria das Filas",
    "6.2. Modelos Markovianos"
)

$target.Range.Text = [string]::Join($lb, $lines)
